# Post data input clean up
# Adds two new analysis columns ("Was the urchin deterred during video?" and
# "Was the urchin sucessfully deterred during video?") to the "Videos" sheet,
# right before the existing "Comments" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Videos")

# Insert two new blank columns at AA (this pushes the old AA "Comments" column to AC)
$ws.Columns("AA:AB").Insert()

# Give the two new columns the same look as the neighboring "Z" (when was kelp
# consumed?) column
$ws.Columns("AA:AB").ColumnWidth = 29.3

# Header row
$ws.Range("AA1").Value = "Was the urchin deterred during video?"
$ws.Range("AB1").Value = "Was the urchin sucessfully deterred during video?"

# Row 2 formulas are standalone (first row of data, entered by hand)
$ws.Range("AA2").Formula = '=IF(T2=0,"urchin never tried",IF(U2>0,"Yes","No"))'
$ws.Range("AB2").Formula = '=IF(AA2="yes",IF(S2>0,"Corynactis was so close","Corynactis was a monster"),IF(T2>0,"Urchin was a beast","Urchin didn''t even try"))'

# Rows 3-37 share one fill-down formula each, matching how the rest of the sheet
# (columns F, G, Z, etc) is built
$ws.Range("AA3:AA37").Formula = '=IF(T3=0,"urchin never tried",IF(U3>0,"Yes","No"))'
$ws.Range("AB3:AB37").Formula = '=IF(AA3="yes",IF(S3>0,"Corynactis was so close","Corynactis was a monster"),IF(T3>0,"Urchin was a beast","Urchin didn''t even try"))'

# Leave the selection on the last touched cell, like the author did after
# filling the new formulas down column AC (old "Comments" column)
[void]$ws.Activate()
[void]$ws.Range("AC36").Select()
